$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." and must become "Version 2.",
# with the text split into the exact run layout the diff calls for:
#   <proofErr spellStart/> "Versi" "on" <proofErr spellEnd/> " 2"
#   <bookmarkStart/><bookmarkEnd/> "."
#
# Range.InsertXML always re-anchors the freshly inserted fragment AFTER any
# left-over ("remainder") text from the run(s) it touched, so each step below
# is chosen/ordered so the remainder + insertion line up in the order we
# actually want in the final document.

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Step 1: split the "Version" run into two runs "Versi" + "on".
$r1 = $d.Range(5, 7)
$r1.InsertXML($pkgHeader + '<w:p><w:r><w:t>on</w:t></w:r></w:p>' + $pkgFooter)

# Step 2: shrink " 1." down to "2" (consumes the trailing suffix "1.",
# leaving the original " " run as a remainder in front of it).
$r2 = $d.Range(8, 10)
$r2.InsertXML($pkgHeader + '<w:p><w:r><w:t>2</w:t></w:r></w:p>' + $pkgFooter)

# Step 3: append the new trailing "." run after the (untouched) bookmark,
# at the very end of the paragraph's text.
$endPos = $d.Content.End - 1
$r3 = $d.Range($endPos, $endPos)
$r3.InsertAfter(".")

# Step 4: merge the leftover " " and "2" runs into a single " 2" run while
# explicitly rewriting the bookmark in place, so it stays correctly
# positioned between " 2" and the trailing "." run.
$r4 = $d.Range(7, 10)
$r4.InsertXML($pkgHeader + '<w:p><w:r><w:t xml:space="preserve"> 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>.</w:t></w:r></w:p>' + $pkgFooter)
